$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update destination names to include accents and extended region info
$ws.Range("B2").Value = "Cancún (y alrededores), México"
$ws.Range("B3").Value = "Ciudad de México (y alrededores), Distrito Federal, México"
$ws.Range("B4").Value = "Bogotá (y alrededores), Colombia"
$ws.Range("B5").Value = "Los Ángeles (y alrededores), California, Estados Unidos de América"
$ws.Range("B6").Value = "Cancún (y alrededores), México"
$ws.Range("B7").Value = "Ciudad de México (y alrededores), Distrito Federal, México"
$ws.Range("B8").Value = "Bogotá (y alrededores), Colombia"
$ws.Range("B9").Value = "Los Ángeles (y alrededores), California, Estados Unidos de América"

# Fix Adults count for rows 6-9 (was 3, should be 2)
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 2

# Widen column B to fit the longer destination text
# (56.17 character-width units; input value compensates for the COM
# interop's internal pixel-based rounding so the saved width matches)
$ws.Columns.Item(2).ColumnWidth = 55.33
